# Update cryptos list: refresh prices and 1h volume percentages,
# and fix a couple of row orderings (RenderToken/BinanceUSD, VeChain/TrustWalletToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E hold non-numeric-looking text (names, URLs, "  +x.xx%  "),
# so a plain .Value assignment keeps them as text.
# Column D holds price text that sometimes LOOKS like a plain number
# (e.g. "0.617", "1.00", "19.20"); assigning those directly would let the
# engine auto-coerce them to numeric cells and silently drop trailing zeros
# (e.g. "1.00" -> 1). Force text via NumberFormat "@" first, then restore the
# cell to the default "Normal" style so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.400.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.065.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.57%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.98"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.82%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +3.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.370.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.068.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.618.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +17.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0815"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0621"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.24%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").Value = "  +4.41%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0962"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +20.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.449.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.87%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0211"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.31%  "
$ws.Range("E48").Value = "  +4.28%  "
$ws.Range("E49").Value = "  +4.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.90%  "
$ws.Range("E51").Value = "  +1.89%  "
